# Scheduled data refresh: recompute currentAveragePrice* / LevePrice* / LeveProfit*
# columns (H:N) per leve row across the crafter Sheets, pulling the latest
# Universalis market-board snapshot. Values are plain computed numbers (no
# formulas in this workbook), so each affected cell is overwritten directly.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# ALC row 15: Morning Glass of Ether / Ether
$ws.Range("H15").Value = 6573.18
$ws.Range("I15").Value = 6573.18
$ws.Range("J15").Value = 0
$ws.Range("K15").Value = 19719.54
$ws.Range("L15").Value = 0
$ws.Range("M15").Value = -19550.54

# ALC row 51: A Bile Business / Shark Oil
$ws.Range("H51").Value = 8347.888999999999
$ws.Range("I51").Value = 3000
$ws.Range("J51").Value = 8662.471
$ws.Range("K51").Value = 3000
$ws.Range("L51").Value = 8662.471
$ws.Range("M51").Value = -2516
$ws.Range("N51").Value = -9630.471

# ALC row 112: Making Ends Meet / Superior Spiritbond Potion
$ws.Range("H112").Value = 4283.148
$ws.Range("I112").Value = 2249
$ws.Range("J112").Value = 4361.385
$ws.Range("K112").Value = 6747
$ws.Range("L112").Value = 13084.155
$ws.Range("M112").Value = -5639
$ws.Range("N112").Value = -15300.155

# ALC row 137: Cutting Edge of Culinary Quality / Magnesia Whetstone
$ws.Range("H137").Value = 3487.0454
$ws.Range("I137").Value = 2580.6875
$ws.Range("J137").Value = 4004.9644
$ws.Range("K137").Value = 7742.0625
$ws.Range("L137").Value = 12014.8932
$ws.Range("M137").Value = -5192.0625
$ws.Range("N137").Value = -17114.8932

# ALC row 138: All-night Crafting / Cunning Craftsman's Tisane
$ws.Range("H138").Value = 4943.757
$ws.Range("I138").Value = 2646.818
$ws.Range("J138").Value = 5344.8096
$ws.Range("K138").Value = 7940.454000000001
$ws.Range("L138").Value = 16034.4288
$ws.Range("M138").Value = -2800.454000000001
$ws.Range("N138").Value = -26314.4288

# ALC row 141: Remedy for Reason / Grade 1 Gemdraught of Mind
$ws.Range("H141").Value = 7275.9287
$ws.Range("I141").Value = 7451
$ws.Range("J141").Value = 5000
$ws.Range("K141").Value = 22353
$ws.Range("L141").Value = 15000
$ws.Range("M141").Value = -17173
$ws.Range("N141").Value = -25360

$ws = $wb.Worksheets.Item("ARM")
# ARM row 32: Ingot We Trust / Steel Ingot
$ws.Range("H32").Value = 3856.5078
$ws.Range("I32").Value = 2150.8035
$ws.Range("J32").Value = 14469.777
$ws.Range("K32").Value = 2150.8035
$ws.Range("L32").Value = 14469.777
$ws.Range("M32").Value = -1863.8035
$ws.Range("N32").Value = -15043.777

# ARM row 132: Don't Bore Me, Ore Me / Mountain Chromite Ingot
$ws.Range("H132").Value = 27850.455
$ws.Range("I132").Value = 31690
$ws.Range("J132").Value = 3533.3333
$ws.Range("K132").Value = 95070
$ws.Range("L132").Value = 10599.9999
$ws.Range("M132").Value = -92540
$ws.Range("N132").Value = -15659.9999

$ws = $wb.Worksheets.Item("BSM")
# BSM row 20: Smelt and Dealt / Iron Ingot
$ws.Range("H20").Value = 2308.7856
$ws.Range("I20").Value = 2049.9473
$ws.Range("J20").Value = 2855.2222
$ws.Range("K20").Value = 2049.9473
$ws.Range("L20").Value = 2855.2222
$ws.Range("M20").Value = -1802.9473
$ws.Range("N20").Value = -3349.2222

# BSM row 95: Crisscrossing / High Steel Kris
$ws.Range("H95").Value = 0
$ws.Range("I95").Value = 0
$ws.Range("J95").Value = 0
$ws.Range("K95").Value = 0
$ws.Range("L95").Value = ""
$ws.Range("N95").Value = 0

# BSM row 134: Ruthenium Supremium / Ruthenium Ingot
$ws.Range("H134").Value = 5245.5835
$ws.Range("I134").Value = 1704.2858
$ws.Range("J134").Value = 6703.7646
$ws.Range("K134").Value = 5112.857400000001
$ws.Range("L134").Value = 20111.2938
$ws.Range("M134").Value = -2577.857400000001
$ws.Range("N134").Value = -25181.2938

$ws = $wb.Worksheets.Item("CRP")
# CRP row 31: Wall Not Found / Walnut Lumber
$ws.Range("H31").Value = 3295.88
$ws.Range("I31").Value = 2491.7407
$ws.Range("J31").Value = 4239.8696
$ws.Range("K31").Value = 2491.7407
$ws.Range("L31").Value = 4239.8696
$ws.Range("M31").Value = -2196.7407
$ws.Range("N31").Value = -4829.8696

# CRP row 34: Armoires of the Rich and Famous / Walnut Lumber
$ws.Range("H34").Value = 3295.88
$ws.Range("I34").Value = 2491.7407
$ws.Range("J34").Value = 4239.8696
$ws.Range("K34").Value = 2491.7407
$ws.Range("L34").Value = 4239.8696
$ws.Range("M34").Value = -2289.7407
$ws.Range("N34").Value = -4643.8696

# CRP row 58: You Do the Heavy Lifting / Mahogany Lumber
$ws.Range("H58").Value = 1430879.9
$ws.Range("I58").Value = 1430879.9
$ws.Range("J58").Value = 0
$ws.Range("K58").Value = 1430879.9
$ws.Range("L58").Value = 0
$ws.Range("M58").Value = ""
$ws.Range("N58").Value = -1430676.9

# CRP row 136: Turali Quality / Dark Mahogany Lumber
$ws.Range("H136").Value = 1430879.9
$ws.Range("I136").Value = 1430879.9
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 4292639.699999999
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = ""
$ws.Range("N136").Value = -4290089.699999999

$ws = $wb.Worksheets.Item("CUL")
# CUL row 68: Such a Butter Face / Fermented Butter
$ws.Range("H68").Value = 974647
$ws.Range("I68").Value = 0
$ws.Range("J68").Value = 974647
$ws.Range("K68").Value = 0
$ws.Range("L68").Value = ""
$ws.Range("M68").Value = 2923941
$ws.Range("N68").Value = -2925563

# CUL row 71: No Margarine of Error (L) / Fermented Butter
$ws.Range("H71").Value = 974647
$ws.Range("I71").Value = 0
$ws.Range("J71").Value = 974647
$ws.Range("K71").Value = 0
$ws.Range("L71").Value = ""
$ws.Range("M71").Value = 8771823
$ws.Range("N71").Value = -8779935

# CUL row 112: Sweet Tooth / Caramels
$ws.Range("H112").Value = 6231.3335
$ws.Range("I112").Value = 5247.5
$ws.Range("J112").Value = 6723.25
$ws.Range("K112").Value = 15742.5
$ws.Range("L112").Value = 20169.75
$ws.Range("M112").Value = -14634.5
$ws.Range("N112").Value = -22385.75

# CUL row 120: A Happy End / Paella
$ws.Range("H120").Value = 13599.857
$ws.Range("I120").Value = 5039.8
$ws.Range("J120").Value = 35000
$ws.Range("K120").Value = 15119.4
$ws.Range("L120").Value = 105000
$ws.Range("M120").Value = -10281.4
$ws.Range("N120").Value = -114676

# CUL row 121: A Cookie for Your Troubles / Coffee Biscuit
$ws.Range("H121").Value = 833
$ws.Range("I121").Value = 0
$ws.Range("J121").Value = 833
$ws.Range("K121").Value = 0
$ws.Range("L121").Value = ""
$ws.Range("M121").Value = 2499
$ws.Range("N121").Value = -5119

# CUL row 133: Friends Are Food / Boiled Alpaca Steak
$ws.Range("H133").Value = 12497.7
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = 12497.7
$ws.Range("K133").Value = 0
$ws.Range("L133").Value = ""
$ws.Range("M133").Value = 37493.10000000001
$ws.Range("N133").Value = -47613.10000000001

$ws = $wb.Worksheets.Item("GSM")
# GSM row 102: Put the Metal to the Peddle / Durium Ingot
$ws.Range("H102").Value = 599874.5
$ws.Range("I102").Value = 1011487.3
$ws.Range("J102").Value = 11856.143
$ws.Range("K102").Value = 1011487.3
$ws.Range("L102").Value = 11856.143
$ws.Range("M102").Value = -1009865.3
$ws.Range("N102").Value = -15100.143

# GSM row 132: On Board for Lar / Lar Ingot
$ws.Range("H132").Value = 3522.6785
$ws.Range("I132").Value = 2895.9443
$ws.Range("J132").Value = 4650.8
$ws.Range("K132").Value = 8687.832900000001
$ws.Range("L132").Value = 13952.4
$ws.Range("M132").Value = -6157.832900000001
$ws.Range("N132").Value = -19012.4

$ws = $wb.Worksheets.Item("LTW")
# LTW row 46: Supply Side Logic / Boar Leather
$ws.Range("H46").Value = 5431.8
$ws.Range("I46").Value = 2031.5
$ws.Range("J46").Value = 6079.476
$ws.Range("K46").Value = 2031.5
$ws.Range("L46").Value = 6079.476
$ws.Range("M46").Value = -1843.5
$ws.Range("N46").Value = -6455.476

# LTW row 132: Tenets of Tanning / Silver Lobo Leather
$ws.Range("H132").Value = 4701.97
$ws.Range("I132").Value = 4229.2163
$ws.Range("J132").Value = 6047.5
$ws.Range("K132").Value = 12687.6489
$ws.Range("L132").Value = 18142.5
$ws.Range("M132").Value = -10157.6489
$ws.Range("N132").Value = -23202.5
